# Apply changes described by the commit:
# "Add SEMICONDUCTOR and PHARMACEUTICAL use type, and modify databases to
#  store the previously hardcoded parameters"
#
# Summary of edits:
#  1. INTERNAL_LOADS sheet: insert two new rows (SEMICONDUCTOR, PHARMACEUTICAL)
#     right after SERVERROOM, copying SERVERROOM's internal-load values, and
#     append five new columns (TCData_Sup, TCData_Re, Hp_X_Cool, Hp_Ratio, N_Hs)
#     for every data row.
#  2. INDOOR_COMFORT sheet: insert the same two new rows, copying SERVERROOM's
#     indoor-comfort values.
#  3. Fix a stray formatting inconsistency on E15 (Qcre_Wm2 for SERVERROOM)
#     so that it carries a border like its column neighbours.
#  4. Re-point the active sheet/selection to match the authored file.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("INTERNAL_LOADS")
$ws2 = $wb.Worksheets.Item("INDOOR_COMFORT")

# ---------------------------------------------------------------------------
# 0. Fix up E15 (SERVERROOM / Qcre_Wm2) so that it has the same bordered
#    style as the rest of column E before we copy row 15's formatting down
#    into the freshly inserted rows 16-17.
# ---------------------------------------------------------------------------
$ws1.Range("E15").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 1. INTERNAL_LOADS: insert two blank rows after SERVERROOM (row 15) and
#    give them the same formatting as row 15 (SERVERROOM).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(16).Resize(2).Insert()
$ws1.Range("A15:N15").Copy()
$ws1.Range("A16:N17").PasteSpecial(-4122)

# Fill in the new SEMICONDUCTOR / PHARMACEUTICAL rows (same internal loads
# as SERVERROOM: only Ed_Wm2=500 and Ea_Wm2=7.1 are non-zero).
$ws1.Range("A16").Value = "SEMICONDUCTOR"
$ws1.Range("A17").Value = "PHARMACEUTICAL"

$newLoadRows = @(16, 17)
foreach ($r in $newLoadRows) {
    $ws1.Cells.Item($r, 2).Value = 0      # Occ_m2p
    $ws1.Cells.Item($r, 3).Value = 0      # Qs_Wp
    $ws1.Cells.Item($r, 4).Value = 0      # X_ghp
    $ws1.Cells.Item($r, 5).Value = 0      # Ea_Wm2
    $ws1.Cells.Item($r, 6).Value = 7.1    # El_Wm2
    $ws1.Cells.Item($r, 7).Value = 0      # Epro_Wm2
    $ws1.Cells.Item($r, 8).Value = 500    # Ed_Wm2
    $ws1.Cells.Item($r, 9).Value = 0      # Vww_ldp
    $ws1.Cells.Item($r, 10).Value = 0     # Vw_ldp
    $ws1.Cells.Item($r, 11).Value = 0     # Qcre_Wm2
    $ws1.Cells.Item($r, 12).Value = 0     # Qhpro_Wm2
    $ws1.Cells.Item($r, 13).Value = 0     # Qcpro_Wm2
    $ws1.Cells.Item($r, 14).Value = 0     # Ev_kWveh
}

# ---------------------------------------------------------------------------
# 2. INTERNAL_LOADS: append the five new columns O:S (TCData_Sup, TCData_Re,
#    Hp_X_Cool, Hp_Ratio, N_Hs) for every code row (2-30).
# ---------------------------------------------------------------------------

# Headers - copy the style of the last existing header (N1) then overwrite text.
$ws1.Range("N1").Copy()
$ws1.Range("O1:S1").PasteSpecial(-4122)
$ws1.Range("O1").Value = "TCData_Sup"
$ws1.Range("P1").Value = "TCData_Re"
$ws1.Range("Q1").Value = "Hp_X_Cool"
$ws1.Range("R1").Value = "Hp_Ratio"
$ws1.Range("S1").Value = "N_Hs"

# Body styles - copy from existing representative columns so fill/border match,
# then apply the correct number formats.
$ws1.Range("B2").Copy()
$ws1.Range("O2:O30").PasteSpecial(-4122)

$ws1.Range("K2").Copy()
$ws1.Range("P2:P30").PasteSpecial(-4122)
$ws1.Range("S2:S30").PasteSpecial(-4122)

$ws1.Range("K2").Copy()
$ws1.Range("Q2:Q30").PasteSpecial(-4122)
$ws1.Range("Q2:Q30").NumberFormat = "0.000"

$ws1.Range("K2").Copy()
$ws1.Range("R2:R30").PasteSpecial(-4122)
$ws1.Range("R2:R30").NumberFormat = "0.00"

# Values - constant for (almost) every row.
$ws1.Range("O2:O30").Value = 7
$ws1.Range("P2:P30").Value = 15
$ws1.Range("Q2:Q30").Value = 0.3
$ws1.Range("R2:R30").Value = 0.83
$ws1.Range("S2:S30").Value = 0.9

# SEMICONDUCTOR / PHARMACEUTICAL have their own Hp_X_Cool values.
$ws1.Range("Q16").Value = 0.222
$ws1.Range("Q17").Value = 0.106

# ---------------------------------------------------------------------------
# 3. INDOOR_COMFORT: insert two blank rows after SERVERROOM (row 15) and
#    give them the same formatting/values as row 15 (SERVERROOM).
# ---------------------------------------------------------------------------
$ws2.Rows.Item(16).Resize(2).Insert()
$ws2.Range("A15:H15").Copy()
$ws2.Range("A16:H17").PasteSpecial(-4122)

$ws2.Range("A16").Value = "SEMICONDUCTOR"
$ws2.Range("A17").Value = "PHARMACEUTICAL"

$newComfortRows = @(16, 17)
foreach ($r in $newComfortRows) {
    $ws2.Cells.Item($r, 2).Value = 24   # Tcs_set_C
    $ws2.Cells.Item($r, 3).Value = 10   # Ths_set_C
    $ws2.Cells.Item($r, 4).Value = 24   # Tcs_setb_C
    $ws2.Cells.Item($r, 5).Value = 10   # Ths_setb_C
    $ws2.Cells.Item($r, 6).Value = 36   # Ve_lsp
    $ws2.Cells.Item($r, 7).Value = 30   # RH_min_pc
    $ws2.Cells.Item($r, 8).Value = 70   # RH_max_pc
}

# ---------------------------------------------------------------------------
# 4. Match the authored file's active sheet / selections.
# ---------------------------------------------------------------------------
$ws2.Range("S23").Select()
$ws1.Activate()
$ws1.Range("X15").Select()
